# Generate Report for Handback
#
# The d62d1220-98b8-4c13-97c3-0e514d99f5a3 file failed its handback transform
# because the handback file name did not match the handoff file name.
# Record the error on the zh-cn and de-de report rows (row 3, which is the
# d62d1220... file) and roll the "Status" for that file up to
# "Handback transform failed" everywhere it is shown (Overview sheet, and
# the per-locale "Status" column).

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: 1zbcjuew.1n4 is different with handoff file name: d62d1220-98b8-4c13-97c3-0e514d99f5a3.70f10829f8a15aacda96055ef5b8a877b5a2307d.zh-cn."
$deError  = "Handback file name: 1zbcjuew.1n4 is different with handoff file name: d62d1220-98b8-4c13-97c3-0e514d99f5a3.70f10829f8a15aacda96055ef5b8a877b5a2307d.de-de."

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Update the status shown on the Overview sheet for the d62d1220... row (row 3).
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Update the per-locale "Status" column (column C) for the same row.
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch in the "Error Detail" column (P).
$wsZh.Range("P3").Value = $zhError
$wsDe.Range("P3").Value = $deError

# Widen the "Error Detail" column so the message is readable.
$wsZh.Range("P:P").ColumnWidth = 39.2
$wsDe.Range("P:P").ColumnWidth = 39.2
